# "Aportes al 6 de enero"
# Applies the Jan 6 2024 (serial 45297) contributions / expenses / receivable
# entries to the workbook, and fixes a mis-attributed "Iverson" entry.

$wb = $excel.ActiveWorkbook

$wsIngreso = $wb.Worksheets.Item("Ingreso")
$wsGastos  = $wb.Worksheets.Item("Gastos")
$wsCxC     = $wb.Worksheets.Item("Cuentas por cobrar")

# ---------------------------------------------------------------------------
# 1) Ingreso!B573 was mis-recorded under the (now retired) "Iverson" member;
#    it should have been credited to "Invitados".
# ---------------------------------------------------------------------------
$wsIngreso.Cells.Item(573, 2).Value = "Invitados"

# ---------------------------------------------------------------------------
# 2) Ingreso: seven new contribution rows dated 2024-01-06 (serial 45297).
# ---------------------------------------------------------------------------

# Rows 582-585 follow the normal "Aporte" pattern (A date / B member /
# C amount / D concepto) -- clone row 203's formatting (unstyled B & C,
# date-styled A, string D) so no new cell styles get minted.
$ingresoRows = @(
    @{ Row = 582; Member = "Rubio"; Amount = 300 },
    @{ Row = 583; Member = "Kawai"; Amount = 100 },
    @{ Row = 584; Member = "Mamao"; Amount = 100 },
    @{ Row = 585; Member = "Punto"; Amount = 100 }
)
foreach ($entry in $ingresoRows) {
    $r = $entry.Row
    $wsIngreso.Range("A203:D203").Copy($wsIngreso.Range("A$r`:D$r"))
    $wsIngreso.Cells.Item($r, 1).Value = 45297
    $wsIngreso.Cells.Item($r, 2).Value = $entry.Member
    $wsIngreso.Cells.Item($r, 3).Value = $entry.Amount
    $wsIngreso.Cells.Item($r, 4).Value = "Aporte"
}

# Rows 586-588 only carry A/B/C (amount 0, no concepto) -- clone the date
# style from row 581's A cell and the amount style from row 1's C cell.
$ingresoZeroRows = @(
    @{ Row = 586; Member = "Javier" },
    @{ Row = 587; Member = "Kibelo" },
    @{ Row = 588; Member = "Joel" }
)
foreach ($entry in $ingresoZeroRows) {
    $r = $entry.Row
    $wsIngreso.Cells.Item(581, 1).Copy($wsIngreso.Cells.Item($r, 1))
    $wsIngreso.Cells.Item($r, 1).Value = 45297
    $wsIngreso.Cells.Item($r, 2).Value = $entry.Member
    $wsIngreso.Cells.Item(1, 3).Copy($wsIngreso.Cells.Item($r, 3))
    $wsIngreso.Cells.Item($r, 3).Value = 0
}

# ---------------------------------------------------------------------------
# 3) Gastos: new expense row 76, 2024-01-06, Arbitro/agua/hielo = 800+150.
# ---------------------------------------------------------------------------
$wsGastos.Range("A75:C75").Copy($wsGastos.Range("A76:C76"))
$wsGastos.Cells.Item(76, 1).Value = 45297
$wsGastos.Cells.Item(76, 2).Value = "Arbitro, agua y hielo"
$wsGastos.Cells.Item(76, 3).Formula = "=800+150"

# ---------------------------------------------------------------------------
# 4) Cuentas por cobrar: new receivable row 7 -- Mamao, Tecnica, 100, with a
#    comment describing the incident.
# ---------------------------------------------------------------------------
$wsCxC.Cells.Item(6, 1).Copy($wsCxC.Cells.Item(7, 1))
$wsCxC.Cells.Item(7, 1).Value = 44932
$wsCxC.Cells.Item(7, 2).Value = "Mamao"
$wsCxC.Cells.Item(7, 3).Value = "Tecnica"
$wsCxC.Cells.Item(7, 4).Value = 100
$wsCxC.Cells.Item(7, 6).Value = "Se sentó en medio del juego"

# ---------------------------------------------------------------------------
# 5) Restore view/selection state per sheet (visiting each sheet to set its
#    own remembered selection), finishing back on the originally-active
#    "Ingreso" tab.
# ---------------------------------------------------------------------------
$wsGastos.Activate()
$wsGastos.Range("A76").Select()

$wsCxC.Activate()
$wsCxC.Range("F7").Select()

$wsIngreso.Activate()
$wsIngreso.Range("D586").Select()

Write-Output "Aportes al 6 de enero applied"
